# CallX/inputs/TestData.xlsx - "done some small changes"
#
# Summary of the edit (reverse engineered from the OOXML diff):
#  - The "Reports" sheet's "Campaigns By Publisher" report rows (rows 2-7)
#    are replaced by "Campaigns" report rows (the same values that already
#    exist in Sheet1 rows 1-6), and the "Campaigns By Publisher" rows are
#    (re)created inside "Sheet1" as new rows 7-12.
#  - Various sheet selections / the active sheet change.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Reports" sheet: rows 2-7 change from "Campaigns By Publisher" data
#    to "Campaigns" data (J also changes 16 -> 15).
# ---------------------------------------------------------------------
$reports = $wb.Worksheets.Item("Reports")
$reports.Activate()

$reportsData = @(
    @("STAGE","Campaigns","Today","Campaigns_Table_Header","Campaigns_Table_Data","campaigns","today","today","campaign_name",15),
    @("STAGE","Campaigns","Yesterday","Campaigns_Table_Header","Campaigns_Table_Data","campaigns","yesterday","yesterday","campaign_name",15),
    @("STAGE","Campaigns","Last_Seven_Days","Campaigns_Table_Header","Campaigns_Table_Data","campaigns","seven-days","yesterday","campaign_name",15),
    @("STAGE","Campaigns","Last_Thirty_Days","Campaigns_Table_Header","Campaigns_Table_Data","campaigns","thirty-days","yesterday","campaign_name",15),
    @("STAGE","Campaigns","Last_Month","Campaigns_Table_Header","Campaigns_Table_Data","campaigns","last-month","last-month_last","campaign_name",15),
    @("STAGE","Campaigns","This_Month","Campaigns_Table_Header","Campaigns_Table_Data","campaigns","this-month","today","campaign_name",15)
)

for ($i = 0; $i -lt $reportsData.Length; $i++) {
    $r = 2 + $i
    $row = $reportsData[$i]
    $reports.Range("A$r").Value = $row[0]
    $reports.Range("B$r").Value = $row[1]
    $reports.Range("C$r").Value = $row[2]
    $reports.Range("D$r").Value = $row[3]
    $reports.Range("E$r").Value = $row[4]
    $reports.Range("F$r").Value = $row[5]
    $reports.Range("G$r").Value = $row[6]
    $reports.Range("H$r").Value = $row[7]
    $reports.Range("I$r").Value = $row[8]
    $reports.Range("J$r").Value = $row[9]
}

# ---------------------------------------------------------------------
# 2) "Sheet1": insert the "Campaigns By Publisher" block as new rows
#    7-12 (rows 13+ already exist and keep their own row numbers).
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()

# Copy the formatting of the row above down into the new rows first so the
# new cells pick up the same style (s="4") as the rest of the table.
$sheet1.Range("A6:J6").Copy()
$sheet1.Range("A7:J12").PasteSpecial(-4122) | Out-Null
$sheet1.Application.CutCopyMode = $false

$sheet1Data = @(
    @("STAGE","Campaigns_By_Publisher","Today","Campaigns_By_Publisher_Table_Header","Campaigns_By_Publisher_Table_Data","campaignbypublisher","today","today","publisher_name",16),
    @("STAGE","Campaigns_By_Publisher","Yesterday","Campaigns_By_Publisher_Table_Header","Campaigns_By_Publisher_Table_Data","campaignbypublisher","yesterday","yesterday","publisher_name",16),
    @("STAGE","Campaigns_By_Publisher","Last_Seven_Days","Campaigns_By_Publisher_Table_Header","Campaigns_By_Publisher_Table_Data","campaignbypublisher","seven-days","yesterday","publisher_name",16),
    @("STAGE","Campaigns_By_Publisher","Last_Thirty_Days","Campaigns_By_Publisher_Table_Header","Campaigns_By_Publisher_Table_Data","campaignbypublisher","thirty-days","yesterday","publisher_name",16),
    @("STAGE","Campaigns_By_Publisher","Last_Month","Campaigns_By_Publisher_Table_Header","Campaigns_By_Publisher_Table_Data","campaignbypublisher","last-month","last-month_last","publisher_name",16),
    @("STAGE","Campaigns_By_Publisher","This_Month","Campaigns_By_Publisher_Table_Header","Campaigns_By_Publisher_Table_Data","campaignbypublisher","this-month","today","publisher_name",16)
)

for ($i = 0; $i -lt $sheet1Data.Length; $i++) {
    $r = 7 + $i
    $row = $sheet1Data[$i]
    $sheet1.Range("A$r").Value = $row[0]
    $sheet1.Range("B$r").Value = $row[1]
    $sheet1.Range("C$r").Value = $row[2]
    $sheet1.Range("D$r").Value = $row[3]
    $sheet1.Range("E$r").Value = $row[4]
    $sheet1.Range("F$r").Value = $row[5]
    $sheet1.Range("G$r").Value = $row[6]
    $sheet1.Range("H$r").Value = $row[7]
    $sheet1.Range("I$r").Value = $row[8]
    $sheet1.Range("J$r").Value = $row[9]
}

# ---------------------------------------------------------------------
# 3) Sheet selections (and active sheet / tab).
# ---------------------------------------------------------------------

# Reports: D2 -> B15
$reports.Activate()
$reports.Range("B15").Select()

# Sales: loses tabSelected, selection H11 -> B2
$sales = $wb.Worksheets.Item("Sales")
$sales.Activate()
$sales.Range("B2").Select()

# Sheet1: topLeftCell cleared, selection D31 -> range A1:J6
$sheet1.Activate()
$sheet1.Range("A1:J6").Select()

# Key Press becomes the active / tabSelected sheet, selection F17 -> J3
$keyPress = $wb.Worksheets.Item("Key Press")
$keyPress.Activate()
$keyPress.Range("J3").Select()
